# BOM.xlsx update — "schematics > rgb led / pcb > via size route"
#
# The piezo buzzer BOM line (row 9) is re-sourced from a PUI Audio
# SMT-1141-T-3-R part to a Mallory Sonalert AST1109MLTRQ part: new
# Mouser link + part number, a lower unit price (3.9 -> 3.12), and the
# row/total/hyperlink bookkeeping that follows from that.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 9 (piezo buzzer) : link (C9) + reference (D9) + price (F9) ---
$ws.Range("C9").Value = "http://www.mouser.ch/ProductDetail/Mallory-Sonalert/AST1109MLTRQ/?qs=sGAEpiMZZMtWZVZ%2fjgUYS%252bu1KhIEHEMRsnSRypyJqVQ%3d"
$ws.Range("D9").Value = "AST1109MLTRQ"
$ws.Range("F9").Value = 3.12

# Row grew very slightly (wrapped link text) in the source edit.
$ws.Rows.Item(9).RowHeight = 15.7

# --- Hyperlinks: rebuild in ref order so the new C9 link lands at rId3
#     and the following links (C11/C12/C14) shift down by one, exactly
#     like the diff shows. ---
$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("C5"), "http://www.mouser.ch/Search/ProductDetail.aspx?qs=8%2f1pEl6ptNseo9Gxrhu%2fPA%3d%3d", "", "", "http://www.mouser.ch/Search/ProductDetail.aspx?qs=8%2f1pEl6ptNseo9Gxrhu%2fPA%3d%3d") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C7"), "http://www.mouser.ch/ProductDetail/Maxim-Integrated/MAX1595EUA33+/?qs=sGAEpiMZZMtitjHzVIkrqUmW7fHvDhXHgnQoEKfsHaU%3d", "", "", "http://www.mouser.ch/ProductDetail/Maxim-Integrated/MAX1595EUA33+/?qs=sGAEpiMZZMtitjHzVIkrqUmW7fHvDhXHgnQoEKfsHaU%3d") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C9"), "http://www.mouser.ch/ProductDetail/Mallory-Sonalert/AST1109MLTRQ/?qs=sGAEpiMZZMtWZVZ%2fjgUYS%252bu1KhIEHEMRsnSRypyJqVQ%3d", "", "", "http://www.mouser.ch/ProductDetail/Mallory-Sonalert/AST1109MLTRQ/?qs=sGAEpiMZZMtWZVZ%2fjgUYS%252bu1KhIEHEMRsnSRypyJqVQ%3d") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C11"), "http://www.mouser.ch/ProductDetail/CK-Components/SK-12C0405-SG-15-RT/?qs=sGAEpiMZZMtHXLepoqNyVaknRufv4Zo6J8yLuspm3Zw%3d", "", "", "http://www.mouser.ch/ProductDetail/CK-Components/SK-12C0405-SG-15-RT/?qs=sGAEpiMZZMtHXLepoqNyVaknRufv4Zo6J8yLuspm3Zw%3d") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C12"), "http://www.mouser.ch/ProductDetail/CK-Components/KMR631NG-ULC-LFS/?qs=sGAEpiMZZMsgGjVA3toVBJ1OkFFtNMGB4KijNZUSro0%3d", "", "", "http://www.mouser.ch/ProductDetail/CK-Components/KMR631NG-ULC-LFS/?qs=sGAEpiMZZMsgGjVA3toVBJ1OkFFtNMGB4KijNZUSro0%3d") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C14"), "http://www.mouser.ch/ProductDetail/Linx-Technologies/BAT-HLD-002-SMT/?qs=%2fha2pyFaduilhNkyJFgy2WekJWQQ7JGY1Lox0Z3adM0%3d", "", "", "http://www.mouser.ch/ProductDetail/Linx-Technologies/BAT-HLD-002-SMT/?qs=%2fha2pyFaduilhNkyJFgy2WekJWQQ7JGY1Lox0Z3adM0%3d") | Out-Null

# --- Column widths (B/C/D in the diff's post-edit numbering == C/D/E
#     here) nudged wider a touch, matching the table's post-edit relayout. ---
$ws.Columns.Item(3).ColumnWidth = 10.729629629629667
$ws.Columns.Item(4).ColumnWidth = 21.21481481481477
$ws.Columns.Item(5).ColumnWidth = 41.207407407407366

# --- View/selection state left at the spot the edit was made from. ---
$ws.Range("F10").Select()

# --- Window tab-ratio tweak (cosmetic UI state). ---
$excel.ActiveWindow.TabRatio = 0.993
